$wb = $excel.ActiveWorkbook

# --- Sheet 1: DQ_Report ---
$ws1 = $wb.Worksheets.Item("DQ_Report")

# Header fix (remove umlaut)
$ws1.Range("B1").Value = "ICD_primaerkode"

# Updated dq_msg text in column D
$ws1.Range("D2").Value = "Kodierung nicht eindeutig. Relation G70 - 586 ist im BfArM nicht vorhanden. "
$ws1.Range("D3").Value = "Kodierung nicht eindeutig. Relation G70 - 589 ist im BfArM nicht vorhanden. "
$ws1.Range("D4").Value = "Kodierung ist nicht eindeutig. Relation E84.80 - 588 ist im BfArM nicht vorhanden.  "
$ws1.Range("D5").Value = "Kodierung ist nicht eindeutig. Relation E75.2 - 325 ist im BfArM nicht vorhanden.  "
$ws1.Range("D6").Value = "Kodierung ist nicht eindeutig. Relation E75.2 - 320 ist im BfArM nicht vorhanden.  "
$ws1.Range("D7").Value = "Fehlendes ICD10 Code.  "
$ws1.Range("D8").Value = "Orpha Kodierung 587 ist im BfArM-Mapping nicht enthalten.  Fehlendes ICD10 Code.  "
$ws1.Range("D9").Value = "Fehlendes Orpha_Kode.  "
$ws1.Range("D10").Value = "Fehlendes Orpha_Kode.  "
$ws1.Range("D11").Value = "Fehlendes ICD10 Code.  "
$ws1.Range("D12").Value = "Kodierung ist nicht eindeutig. Relation E66.89 - 320 ist im BfArM nicht vorhanden.  "
$ws1.Range("D13").Value = "Fehlendes Orpha_Kode.  "
$ws1.Range("D14").Value = "Fehlendes Orpha_Kode.  "
$ws1.Range("D15").Value = "Kodierung ist nicht eindeutig. Relation E85.0 - 586 ist im BfArM nicht vorhanden.  "

# Remove the now-obsolete last row (row 16)
$ws1.Rows.Item(16).Delete()

# --- Sheet 2: Statistik ---
$ws2 = $wb.Worksheets.Item("Statistik")

$ws2.Range("F1").Value = "K2_icdRd_no"
$ws2.Range("G1").Value = "K3_rd_no"

$ws2.Range("E2").Value = 97.7
$ws2.Range("G2").Value = 297
